# Update vaccine dose by age: append two new weekly rows
# (2021-09-06 / serial 44445 and 2021-09-13 / serial 44452)
# to both the "1st dose" and "2nd dose" sheets, matching the formatting
# already used by the rest of the table (copy formats down from the last
# "plain" data row so no new styles are introduced).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "1st dose" -> new rows 20 and 21 (plain values)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1st dose")

$ws1.Range("A17:H17").Copy()
$ws1.Range("A20:H20").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("A17:H17").Copy()
$ws1.Range("A21:H21").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws1.Range("A20").Value = 44445
$ws1.Range("B20").Value = 3497433
$ws1.Range("C20").Value = 3091102
$ws1.Range("D20").Value = 3563988
$ws1.Range("E20").Value = 7843679
$ws1.Range("F20").Value = 6656881
$ws1.Range("G20").Value = 3474049
$ws1.Range("H20").Value = 1863037

$ws1.Range("A21").Value = 44452
$ws1.Range("B21").Value = 4459494
$ws1.Range("C21").Value = 3912573
$ws1.Range("D21").Value = 4845871
$ws1.Range("E21").Value = 7895083
$ws1.Range("F21").Value = 6677000
$ws1.Range("G21").Value = 3480174
$ws1.Range("H21").Value = 1867227

# ---------------------------------------------------------------------------
# Sheet "2nd dose" -> new rows 20 and 21 (formulas: reported total - correction)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2nd dose")

$ws2.Range("A17:H17").Copy()
$ws2.Range("A20:H20").PasteSpecial(-4122) # xlPasteFormats
$ws2.Range("A17:H17").Copy()
$ws2.Range("A21:H21").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws2.Range("A20").Value = 44445
$ws2.Range("B20").Formula = "=1584973-43"
$ws2.Range("C20").Formula = "=1772048-860571"
$ws2.Range("D20").Formula = "=1709079-251759"
$ws2.Range("E20").Formula = "=1542659-80621"
$ws2.Range("F20").Formula = "=6040407-65402"
$ws2.Range("G20").Formula = "=3309226-9060"
$ws2.Range("H20").Formula = "=1780185-1134"

$ws2.Range("A21").Value = 44452
$ws2.Range("B21").Formula = "=1912209-66"
$ws2.Range("C21").Formula = "=2011816-887563"
$ws2.Range("D21").Formula = "=2014986-276082"
$ws2.Range("E21").Formula = "=2814585-93245"
$ws2.Range("F21").Formula = "=6167274-68803"
$ws2.Range("G21").Formula = "=3338637-9616"
$ws2.Range("H21").Formula = "=1787410-1486"

# Leave the cursor where the author would naturally land after typing the
# last new row on each sheet (cosmetic only - does not affect data/styles,
# and does not disturb which sheet/tab is active).
$ws1.Range("F21").Select()
$ws2.Range("D15").Select()

$wb.Save()
